# Auto-generated Excel COM-interop script applying cryptos.xlsx price/volume update
# (commit: "Updated cryptos list on Thu Dec 14 03:37:37 UTC 2023 with GitHub Actions")

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "'42.741.64"
$ws.Cells.Item(2, 4).Style = "Normal"
$ws.Cells.Item(2, 5).Value = "'  +4.41%  "
$ws.Cells.Item(2, 5).Style = "Normal"
$ws.Cells.Item(3, 4).Value = "'2.252.59"
$ws.Cells.Item(3, 4).Style = "Normal"
$ws.Cells.Item(3, 5).Value = "'  +3.91%  "
$ws.Cells.Item(3, 5).Style = "Normal"
$ws.Cells.Item(4, 5).Value = "'  -0.02%  "
$ws.Cells.Item(4, 5).Style = "Normal"
$ws.Cells.Item(5, 4).Value = "'249.01"
$ws.Cells.Item(5, 4).Style = "Normal"
$ws.Cells.Item(5, 5).Value = "'  +0.48%  "
$ws.Cells.Item(5, 5).Style = "Normal"
$ws.Cells.Item(6, 4).Value = "'0.622"
$ws.Cells.Item(6, 4).Style = "Normal"
$ws.Cells.Item(6, 5).Value = "'  +0.92%  "
$ws.Cells.Item(6, 5).Style = "Normal"
$ws.Cells.Item(7, 4).Value = "'69.95"
$ws.Cells.Item(7, 4).Style = "Normal"
$ws.Cells.Item(7, 5).Value = "'  +5.36%  "
$ws.Cells.Item(7, 5).Style = "Normal"
$ws.Cells.Item(8, 5).Value = "'  -0.15%  "
$ws.Cells.Item(8, 5).Style = "Normal"
$ws.Cells.Item(9, 4).Value = "'0.656"
$ws.Cells.Item(9, 4).Style = "Normal"
$ws.Cells.Item(9, 5).Value = "'  +17.02%  "
$ws.Cells.Item(9, 5).Style = "Normal"
$ws.Cells.Item(10, 4).Value = "'39.10"
$ws.Cells.Item(10, 4).Style = "Normal"
$ws.Cells.Item(10, 5).Value = "'  +10.79%  "
$ws.Cells.Item(10, 5).Style = "Normal"
$ws.Cells.Item(11, 2).Value = "OKB"
$ws.Cells.Item(11, 3).Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Cells.Item(11, 4).Value = "'59.23"
$ws.Cells.Item(11, 4).Style = "Normal"
$ws.Cells.Item(11, 5).Value = "'  +2.48%  "
$ws.Cells.Item(11, 5).Style = "Normal"
$ws.Cells.Item(12, 2).Value = "Dogecoin"
$ws.Cells.Item(12, 3).Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Cells.Item(12, 4).Value = "'0.0964"
$ws.Cells.Item(12, 4).Style = "Normal"
$ws.Cells.Item(12, 5).Value = "'  +4.66%  "
$ws.Cells.Item(12, 5).Style = "Normal"
$ws.Cells.Item(13, 4).Value = "'7.47"
$ws.Cells.Item(13, 4).Style = "Normal"
$ws.Cells.Item(13, 5).Value = "'  +9.00%  "
$ws.Cells.Item(13, 5).Style = "Normal"
$ws.Cells.Item(14, 5).Value = "'  +0.41%  "
$ws.Cells.Item(14, 5).Style = "Normal"
$ws.Cells.Item(15, 4).Value = "'2.580.80"
$ws.Cells.Item(15, 4).Style = "Normal"
$ws.Cells.Item(15, 5).Value = "'  +3.50%  "
$ws.Cells.Item(15, 5).Style = "Normal"
$ws.Cells.Item(16, 4).Value = "'14.72"
$ws.Cells.Item(16, 4).Style = "Normal"
$ws.Cells.Item(16, 5).Value = "'  +3.98%  "
$ws.Cells.Item(16, 5).Style = "Normal"
$ws.Cells.Item(17, 5).Value = "'  +3.16%  "
$ws.Cells.Item(17, 5).Style = "Normal"
$ws.Cells.Item(18, 4).Value = "'2.258.49"
$ws.Cells.Item(18, 4).Style = "Normal"
$ws.Cells.Item(18, 5).Value = "'  +3.50%  "
$ws.Cells.Item(18, 5).Style = "Normal"
$ws.Cells.Item(19, 4).Value = "'42.652.55"
$ws.Cells.Item(19, 4).Style = "Normal"
$ws.Cells.Item(19, 5).Value = "'  +4.40%  "
$ws.Cells.Item(19, 5).Style = "Normal"
$ws.Cells.Item(20, 5).Value = "'  +5.17%  "
$ws.Cells.Item(20, 5).Style = "Normal"
$ws.Cells.Item(21, 4).Value = "'6.28"
$ws.Cells.Item(21, 4).Style = "Normal"
$ws.Cells.Item(21, 5).Value = "'  +3.63%  "
$ws.Cells.Item(21, 5).Style = "Normal"
$ws.Cells.Item(22, 4).Value = "'72.80"
$ws.Cells.Item(22, 4).Style = "Normal"
$ws.Cells.Item(22, 5).Value = "'  +2.03%  "
$ws.Cells.Item(22, 5).Style = "Normal"
$ws.Cells.Item(23, 4).Value = "'232.92"
$ws.Cells.Item(23, 4).Style = "Normal"
$ws.Cells.Item(23, 5).Value = "'  +1.72%  "
$ws.Cells.Item(23, 5).Style = "Normal"
$ws.Cells.Item(24, 5).Value = "'  +1.37%  "
$ws.Cells.Item(24, 5).Style = "Normal"
$ws.Cells.Item(25, 4).Value = "'3.94"
$ws.Cells.Item(25, 4).Style = "Normal"
$ws.Cells.Item(25, 5).Value = "'  +5.37%  "
$ws.Cells.Item(25, 5).Style = "Normal"
$ws.Cells.Item(26, 4).Value = "'11.53"
$ws.Cells.Item(26, 4).Style = "Normal"
$ws.Cells.Item(26, 5).Value = "'  -0.33%  "
$ws.Cells.Item(26, 5).Style = "Normal"
$ws.Cells.Item(27, 5).Value = "'  +0.11%  "
$ws.Cells.Item(27, 5).Style = "Normal"
$ws.Cells.Item(28, 5).Value = "'  +0.39%  "
$ws.Cells.Item(28, 5).Style = "Normal"
$ws.Cells.Item(29, 4).Value = "'3.64"
$ws.Cells.Item(29, 4).Style = "Normal"
$ws.Cells.Item(29, 5).Value = "'  -1.81%  "
$ws.Cells.Item(29, 5).Style = "Normal"
$ws.Cells.Item(30, 4).Value = "'2.11"
$ws.Cells.Item(30, 4).Style = "Normal"
$ws.Cells.Item(30, 5).Value = "'  -0.64%  "
$ws.Cells.Item(30, 5).Style = "Normal"
$ws.Cells.Item(31, 4).Value = "'166.91"
$ws.Cells.Item(31, 4).Style = "Normal"
$ws.Cells.Item(31, 5).Value = "'  -0.77%  "
$ws.Cells.Item(31, 5).Style = "Normal"
$ws.Cells.Item(32, 4).Value = "'20.91"
$ws.Cells.Item(32, 4).Style = "Normal"
$ws.Cells.Item(32, 5).Value = "'  +3.88%  "
$ws.Cells.Item(32, 5).Style = "Normal"
$ws.Cells.Item(33, 4).Value = "'6.36"
$ws.Cells.Item(33, 4).Style = "Normal"
$ws.Cells.Item(33, 5).Value = "'  +14.25%  "
$ws.Cells.Item(33, 5).Style = "Normal"
$ws.Cells.Item(34, 5).Value = "'  +5.68%  "
$ws.Cells.Item(34, 5).Style = "Normal"
$ws.Cells.Item(35, 4).Value = "'31.43"
$ws.Cells.Item(35, 4).Style = "Normal"
$ws.Cells.Item(35, 5).Value = "'  +23.34%  "
$ws.Cells.Item(35, 5).Style = "Normal"
$ws.Cells.Item(36, 4).Value = "'0.0790"
$ws.Cells.Item(36, 4).Style = "Normal"
$ws.Cells.Item(36, 5).Value = "'  +7.76%  "
$ws.Cells.Item(36, 5).Style = "Normal"
$ws.Cells.Item(37, 5).Value = "'  +4.14%  "
$ws.Cells.Item(37, 5).Style = "Normal"
$ws.Cells.Item(38, 4).Value = "'4.39"
$ws.Cells.Item(38, 4).Style = "Normal"
$ws.Cells.Item(38, 5).Value = "'  +8.91%  "
$ws.Cells.Item(38, 5).Style = "Normal"
$ws.Cells.Item(40, 5).Value = "'  +6.59%  "
$ws.Cells.Item(40, 5).Style = "Normal"
$ws.Cells.Item(41, 4).Value = "'2.31"
$ws.Cells.Item(41, 4).Style = "Normal"
$ws.Cells.Item(41, 5).Value = "'  +6.97%  "
$ws.Cells.Item(41, 5).Style = "Normal"
$ws.Cells.Item(42, 4).Value = "'12.49"
$ws.Cells.Item(42, 4).Style = "Normal"
$ws.Cells.Item(42, 5).Value = "'  +8.91%  "
$ws.Cells.Item(42, 5).Style = "Normal"
$ws.Cells.Item(44, 4).Value = "'62.29"
$ws.Cells.Item(44, 4).Style = "Normal"
$ws.Cells.Item(44, 5).Value = "'  +4.09%  "
$ws.Cells.Item(44, 5).Style = "Normal"
$ws.Cells.Item(45, 4).Value = "'9.07"
$ws.Cells.Item(45, 4).Style = "Normal"
$ws.Cells.Item(45, 5).Value = "'  +7.07%  "
$ws.Cells.Item(45, 5).Style = "Normal"
$ws.Cells.Item(46, 2).Value = "Algorand"
$ws.Cells.Item(46, 3).Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Cells.Item(46, 4).Value = "'0.200"
$ws.Cells.Item(46, 4).Style = "Normal"
$ws.Cells.Item(46, 5).Value = "'  +5.55%  "
$ws.Cells.Item(46, 5).Style = "Normal"
$ws.Cells.Item(47, 2).Value = "FTXToken"
$ws.Cells.Item(47, 3).Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Cells.Item(47, 4).Value = "'4.84"
$ws.Cells.Item(47, 4).Style = "Normal"
$ws.Cells.Item(47, 5).Value = "'  +2.05%  "
$ws.Cells.Item(47, 5).Style = "Normal"
$ws.Cells.Item(48, 5).Value = "'  +4.07%  "
$ws.Cells.Item(48, 5).Style = "Normal"
$ws.Cells.Item(49, 4).Value = "'1.00"
$ws.Cells.Item(49, 4).Style = "Normal"
$ws.Cells.Item(49, 5).Value = "'  -0.41%  "
$ws.Cells.Item(49, 5).Style = "Normal"
$ws.Cells.Item(50, 4).Value = "'1.16"
$ws.Cells.Item(50, 4).Style = "Normal"
$ws.Cells.Item(50, 5).Value = "'  +1.00%  "
$ws.Cells.Item(50, 5).Style = "Normal"
$ws.Cells.Item(51, 5).Value = "'  +3.87%  "
$ws.Cells.Item(51, 5).Style = "Normal"
